$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B3").ClearContents()
$ws.Range("B2:B3").Font.Underline = $true

$ws.Range("B2:B3").Select()
